{"js": "// Remove the trailing \"Ver no Jupiter...\" and \"\u00a9 2020 ...\" footer\n// paragraphs (and the blank paragraph that precedes them), which were\n// dropped from the page when the site was rebuilt.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst TARGET_TEXTS = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\nlet removeIndices = [];\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text.trim();\n  if (TARGET_TEXTS.includes(text)) {\n    removeIndices.push(i);\n  }\n}\n\n// Also remove the blank paragraph immediately before the first matched\n// paragraph (it was the spacer between \"M\u00c1QUINASEscola PRO-TEC\" and the\n// footer block, and the diff drops it together with the footer text).\nif (removeIndices.length > 0) {\n  const firstIdx = removeIndices[0];\n  const prev = items[firstIdx - 1];\n  if (prev && prev.text.trim() === \"\") {\n    removeIndices.unshift(firstIdx - 1);\n  }\n}\n\nfor (const idx of removeIndices) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" and \"(c) 2020 ...\" footer\n# paragraphs (plus the blank spacer paragraph right before them), which\n# were dropped from the page when the site was rebuilt.\n$d = $word.ActiveDocument\n\n$paras = @($d.Paragraphs)\n$removeIndices = New-Object System.Collections.ArrayList\n\nfor ($i = 0; $i -lt $paras.Count; $i++) {\n    $text = $paras[$i].Range.Text\n    if ($text -like \"*Ver no Jupiter Salvar em pdf Salvar em docx*\") {\n        [void]$removeIndices.Add($i)\n    } elseif ($text -like \"*Contact: luizeleno@usp.br*\") {\n        [void]$removeIndices.Add($i)\n    }\n}\n\n# Also drop the blank paragraph right before the first matched paragraph\n# (the spacer between \"MAQUINASEscola PRO-TEC\" and the footer block).\nif ($removeIndices.Count -gt 0) {\n    $firstIndex = $removeIndices[0]\n    if ($firstIndex -gt 0 -and $paras[$firstIndex - 1].Range.Text.Trim().Length -eq 0) {\n        [void]$removeIndices.Insert(0, $firstIndex - 1)\n    }\n}\n\n# Delete from the highest index down so earlier (lower) indices remain valid.\n$sorted = $removeIndices | Sort-Object -Descending\nforeach ($idx in $sorted) {\n    $paras[$idx].Range.Delete()\n}\n"}
